$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.877819061279297
$ws.Range("C2").Value = 5.522988319396973
$ws.Range("D2").Value = 13.01597785949707
$ws.Range("E2").Value = 57.85714340209961
